# Insert one new data row into the "Zanahoria" (Carrot) price sheet.
#
# The new record (week of 2023-12-15, serial 45275) is inserted as row 534,
# pushing the previously existing rows 534-651 down to rows 535-652.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 534; existing rows 534.. shift down.
$ws.Rows(534).Insert()

# Populate the newly inserted row 534 with the new data.
$ws.Cells.Item(534, 1).Value = 5
$ws.Cells.Item(534, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(534, 3).Value = "Maule"
$ws.Cells.Item(534, 4).Value = 45275
$ws.Cells.Item(534, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(534, 5).Value = 7
$ws.Cells.Item(534, 6).Value = 100114013
$ws.Cells.Item(534, 7).Value = "Zanahoria"
$ws.Cells.Item(534, 8).Value = "Sin especificar"
$ws.Cells.Item(534, 9).Value = "Primera"
$ws.Cells.Item(534, 10).Value = 700
$ws.Cells.Item(534, 11).Value = 5500
$ws.Cells.Item(534, 12).Value = 5500
$ws.Cells.Item(534, 13).Value = 5500
$ws.Cells.Item(534, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(534, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(534, 16).Value = 275
$ws.Cells.Item(534, 17).Value = 20
$ws.Cells.Item(534, 18).Value = "Hortaliza"
